# Auto update: 2025-12-06 01:15:18
# Update K (최종점수) and N (MACRO_SCORE) columns for rows 2-7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 62.4
    3 = 54.2
    4 = 50.4
    5 = 47.6
    6 = 39.6
    7 = 39.6
}

$macroScore = 50.60178744571824

foreach ($row in 2..7) {
    $ws.Range("K$row").Value = $kValues[$row]
    $ws.Range("N$row").Value = $macroScore
}
